# iron_native PowerPoint COM-interop edit script
#
# Implements (see commit "18/02/25 practica 3 actualizada del 2"):
#  1. Fill in slide 2 ("Mi presentacion") title + body text.
#  2. Add a new slide 3 ("Cosas que me gustan") after slide 2, reusing
#     slide 2's "Title and Content" layout, with the body placeholder
#     repositioned/resized and filled in.
#  3. Refresh the cached datetimeFigureOut date text (13/02/2025 ->
#     18/02/2025) on the slide master and every slide layout.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Slide 2: "Mi presentacion"
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

$title2 = $s2.Shapes.Item(1).TextFrame.TextRange
$title2.Text = "Mi presentacion"
$title2.ParagraphFormat.Alignment = 2   # ppAlignCenter

$body2 = $s2.Shapes.Item(2).TextFrame.TextRange
$body2.Text = "Me llamo Alan Eduardo Leon Merino vivo en Culiacan, Sinaloa pero tengo a toda mi familia  en Oaxaca y me gusta estar allá por que la paso con mi familia y me divierto cuando la paso allá tengo 15 años y todavía no se usar una computadora correctamente y me gustaría aprender mas sobre la programación."
$body2.ParagraphFormat.Alignment = 2    # ppAlignCenter
$body2.Font.Size = 28

# ---------------------------------------------------------------------
# 2. New slide 3: "Cosas que me gustan" (same "Title and Content" layout
#    as slide 2), inserted right after slide 2.
# ---------------------------------------------------------------------
$layout = $s2.CustomLayout
$s3 = $p.Slides.AddSlide(3, $layout)

$title3 = $s3.Shapes.Item(1).TextFrame.TextRange
$title3.Text = " Cosas que me gustan"
$title3.ParagraphFormat.Alignment = 2   # ppAlignCenter

$contentShape = $s3.Shapes.Item(2)
$contentShape.Left = 352697 / 12700
$contentShape.Top = 627018 / 12700
$contentShape.Width = 11403873 / 12700
$contentShape.Height = 5621382 / 12700

$body3 = $contentShape.TextFrame.TextRange
$body3.Text = "`r Me gusta la comida como el pozole , las tlayudas, los huaraches de Oaxaca.`rMe gusta dormir o descansar después de una actividad física y me toma mucho tiempo.`rMe gusta las películas de acción o de suspenso.`rNo me gusta tanto leer, pero si leo me gustan de finanzas o de terror.`r"
$body3.Font.Size = 32

# First paragraph is blank / has no bullet (matches the "No me gusta..."
# body's lead-in blank line in the source deck).
$firstPara3 = $body3.Paragraphs(1, 1)
$firstPara3.ParagraphFormat.Bullet.Type = 0   # ppBulletNone
$ruler3 = $contentShape.TextFrame.Ruler
$lvl3 = $ruler3.Levels.Item(1)
$lvl3.FirstMargin = 0
$lvl3.LeftMargin = 0

# ---------------------------------------------------------------------
# 3. Refresh cached date fields across the slide master and all layouts
#    (datetimeFigureOut placeholders recache the "today" string whenever
#    the deck is touched/saved in PowerPoint).
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes, [string]$newDate) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDate = $false
        try { $isDate = ($shp.PlaceholderFormat.Type -eq 16) } catch {}
        if ($isDate -and $shp.HasTextFrame) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes "18/02/2025"

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $lay = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $lay.Shapes "18/02/2025"
}
